# Update RedDot System and Add ModelManager
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing DateCardsDB row: datedata.xlsx -> dateData.xlsx
$ws.Range("E4").Value = "schedule/dateData.xlsx"

# Add new row for MessageReaderDB / MessageDB (ModelManager entry)
$ws.Range("B5").Value = "MessageReaderDB"
$ws.Range("C5").Value = "MessageDB"
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = "phone/phoneMessageData.xlsx"

# Widen column E to fit the new content
$ws.Columns.Item(5).ColumnWidth = 30.5714285714286

# Move selection to E7, matching the workbook's last saved cursor position
[void]$ws.Range("E7").Select()
